$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update item name and unit text (shared strings "Telur" -> "Air Mineral", "pcs" -> "ml")
$ws.Range("A2").Value = "Air Mineral"
$ws.Range("D2").Value = "ml"

# Update purchase date
$ws.Range("B2").Value = (Get-Date -Year 2023 -Month 3 -Day 17 -Hour 0 -Minute 0 -Second 0).Date

# Update quantity
$ws.Range("C2").Value = 1000

# Update cost
$ws.Range("E2").Value = 3000

# Replace formula in F2 with a plain static value
$ws.Range("F2").Value = 3

# Update selection to F2
$ws.Range("F2").Select()
